# google__python-fire.xlsx : "adding all pytype stubs and extra row to resultant output"
#
# This script:
#  1. Updates the E/F (Scalpel Type / Status) columns for the rows whose
#     inferred-type stub sets gained the missing "real" PyType entries
#     (rows 6/7, 22/23, 50/51) - including flipping the Loss/Win/Neutral
#     verdict + cell color where the extra stub changed the outcome.
#  2. Updates the running tally row (271): PyType Wins 43->42, Scalpel Wins 11->13.
#  3. Repurposes row 272 to show the new "Scalpel Accuracy:" figure and
#     pushes the "Accuracy over PyType" summary down into a new row 273
#     with the updated percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# OLE (BGR) colors matching the existing fill palette used by this sheet:
#   Neutral -> orange (00FFA500)
#   Win     -> green  (00008000)
#   Loss    -> red    (00FF0000)
$colorNeutral = 42495
$colorWin     = 32768
$colorLoss    = 255

# --- Row 6/7 : decorators.py / GetMetadata -> dict -------------------------
$ws.Range("E6").Value = "{'Dict[any, any]', 'any'}"
$ws.Range("E7").Value = "Dict[any, any]"
$ws.Range("F7").Value = "Neutral"
$ws.Range("F7").Interior.Color = $colorNeutral

# --- Row 22/23 : docstrings.py / _join_lines -> Optional -------------------
$ws.Range("E22").Value = "{'str', 'any'}"
$ws.Range("E23").Value = "str"

# --- Row 50/51 : docstrings.py / _get_after_directive -> Any ---------------
$ws.Range("E50").Value = "{'str', 'any'}"
$ws.Range("F50").Value = "Win"
$ws.Range("F50").Interior.Color = $colorWin

$ws.Range("E51").Value = "str"
$ws.Range("F51").Value = "Win"
$ws.Range("F51").Interior.Color = $colorWin

# --- Row 271 : tally row -----------------------------------------------
$ws.Range("D271").Value = 42
$ws.Range("F271").Value = 13

# --- Row 272 : was the "Accuracy over PyType" row, now becomes the new
#     "Scalpel Accuracy:" row, with the accuracy-over-pytype figure moving
#     down to a brand-new row 273. Clone formatting from row 271 first so
#     the new row carries the same cell styling. ---------------------------
$ws.Range("A271:F271").Copy()
$ws.Range("A273:F273").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C272").Value = "Scalpel Accuracy:"
$ws.Range("D272").Value = 540.48
$ws.Range("E272").Value = ""
$ws.Range("F272").Value = ""

$ws.Range("A273").Value = ""
$ws.Range("B273").Value = ""
$ws.Range("C273").Value = ""
$ws.Range("D273").Value = ""
$ws.Range("E273").Value = "Accuracy over PyType"
$ws.Range("F273").Value = 30.95
